$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "Ericsson India Global Services Pvt. Ltd. as a Senior Solution Integrator."
$find.Replacement.Text = "Accenture India Pvt. Ltd. as a Team Lead & MSO."
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
